$d = $word.ActiveDocument

# Map of old "a×b=" expressions to new ones, per the commit diff.
# All 25 values across the document are unique, and no replacement's
# new text collides with another entry's old text, so simple sequential
# Find/Replace (MatchWholeWord) is safe and order-independent.

$replacements = @(
    @("388×3=", "455×7="),
    @("881×7=", "109×3="),
    @("826×3=", "859×2="),
    @("506×2=", "395×9="),
    @("341×5=", "976×2="),
    @("517×6=", "599×3="),
    @("784×5=", "545×9="),
    @("219×6=", "396×9="),
    @("238×2=", "920×6="),
    @("856×8=", "764×3="),
    @("988×2=", "361×2="),
    @("772×8=", "266×3="),
    @("646×7=", "347×3="),
    @("486×2=", "233×2="),
    @("454×8=", "213×8="),
    @("339×2=", "899×3="),
    @("275×3=", "631×3="),
    @("854×8=", "811×3="),
    @("580×7=", "537×3="),
    @("586×8=", "565×8="),
    @("124×8=", "827×9="),
    @("879×3=", "322×7="),
    @("536×3=", "407×7="),
    @("128×5=", "426×7="),
    @("433×7=", "747×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replace '$old' -> '$new': found=$found"
}
